$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.660.59'
$ws.Range("E2").Value = '  -2.19%  '
$ws.Range("D3").Value = '1.760.42'
$ws.Range("E3").Value = '  -3.07%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '323.93'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.33%  '
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4306'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.29%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3603'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.75%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07572'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.47%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.19'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.111'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.62%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.10%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.75'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -6.41%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.072'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.74%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.239'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.98%  '
$ws.Range("D16").Value = '1.759.68'
$ws.Range("E16").Value = '  -3.97%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.39'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.98%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001066'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.48%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06425'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.76%  '
$ws.Range("E20").Value = '  +0.07%  '
$ws.Range("E21").Value = '  -2.34%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.878'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.17%  '
$ws.Range("D23").Value = '27.706.76'
$ws.Range("E23").Value = '  -2.09%  '
$ws.Range("E24").Value = '  -3.30%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.084'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.46%  '
$ws.Range("E26").Value = '  +0.23%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.56'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.93%  '
$ws.Range("D28").Value = '1.959.87'
$ws.Range("E28").Value = '  -3.63%  '
$ws.Range("E29").Value = '  -6.52%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.74'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.34%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.096'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -9.89%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.688'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.75%  '
$ws.Range("E33").Value = '  -6.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08954'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.63%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.20'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02301'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2116'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.96%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06012'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.17%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6354'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.37%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.952'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.76%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.187'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.54%  '
$ws.Range("E42").Value = '  +0.15%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.396'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.38%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.899'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.92%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.32'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.24%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5930'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.21%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.712'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.987'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.77%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '122.76'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.36%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.172'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.18%  '
$ws.Range("E51").Value = '  -1.88%  '
